$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of H1 onto I1:J1 so they match the other headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-43
$iValues = @(5,8,9,6,9,8,7,7,9,9,9,7,6,7,7,9,9,8,10,7,9,6,8,9,8,9,9,8,9,8,8,7,8,7,7,8,8,7,7,8,8,7)
$jValues = @(5,8,9,6,9,8,7,7,9,9,9,7,7,7,8,9,9,8,10,7,9,6,8,9,8,9,9,8,9,8,9,7,8,7,7,8,8,7,7,8,8,7)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
